$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 4.2
$ws.Range("T2").Value = 1.79
$ws.Range("H4").Value = 3.65
$ws.Range("I4").Value = 4.4
$ws.Range("K4").Value = 4.3
$ws.Range("L4").Value = 1.23
$ws.Range("T4").Value = 1.54
$ws.Range("U4").Value = 2.08
$ws.Range("W4").Value = 1.9
$ws.Range("AA4").Value = 90
$ws.Range("F5").Value = 1.9
$ws.Range("H5").Value = 3.85
$ws.Range("I5").Value = 4.7
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 4.2
$ws.Range("N5").Value = 3.15
$ws.Range("V5").Value = 1.27
$ws.Range("H6").Value = 2.96
$ws.Range("K6").Value = 4
$ws.Range("Q6").Value = 1.33
$ws.Range("H7").Value = 1.04
$ws.Range("M7").Value = 1.06
$ws.Range("L11").Value = 1.35
$ws.Range("V11").Value = 2
$ws.Range("T12").Value = 1.56
$ws.Range("U12").Value = 2.46
$ws.Range("F13").Value = 2.08
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 4.1
$ws.Range("J13").Value = 3.55
$ws.Range("Q13").Value = 1.73
$ws.Range("T13").Value = 1.64
$ws.Range("V13").Value = 1.32
$ws.Range("H14").Value = 1.94
$ws.Range("J14").Value = 2.7
$ws.Range("N15").Value = 2.56
$ws.Range("P15").Value = 1.63
$ws.Range("Q15").Value = 1.95
$ws.Range("U15").Value = 1.04
$ws.Range("T16").Value = 1.04
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 2.68
$ws.Range("O17").Value = 1.06
$ws.Range("T17").Value = 1.04
$ws.Range("U17").Value = 1.04
$ws.Range("N18").Value = 2.44
$ws.Range("O18").Value = 1.21
$ws.Range("Q18").Value = 1.38
$ws.Range("R18").Value = 1.39
$ws.Range("S18").Value = 2.1
$ws.Range("T18").Value = 1.04
$ws.Range("U18").Value = 1.04
$ws.Range("G19").Value = 1.37
$ws.Range("J19").Value = 5.5
$ws.Range("F20").Value = 2.78
$ws.Range("H20").Value = 2.56
$ws.Range("P20").Value = 1.73
$ws.Range("Q20").Value = 2.08
$ws.Range("T20").Value = 1.04
$ws.Range("U20").Value = 1.04
$ws.Range("M21").Value = 1.08
$ws.Range("N21").Value = 2.42
$ws.Range("O21").Value = 1.35
$ws.Range("T21").Value = 1.04
$ws.Range("U21").Value = 1.04
$ws.Range("H22").Value = 1.74
$ws.Range("K22").Value = 6.2
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 2.64
$ws.Range("O22").Value = 1.06
$ws.Range("T22").Value = 1.04
$ws.Range("U22").Value = 1.73
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 2.5
$ws.Range("O23").Value = 1.06
$ws.Range("Q23").Value = 1.76
$ws.Range("T23").Value = 1.04
$ws.Range("U23").Value = 1.04
$ws.Range("G25").Value = 3.05
$ws.Range("I25").Value = 3.45
$ws.Range("S25").Value = 4.2
$ws.Range("T25").Value = 1.77
$ws.Range("V25").Value = 1.4
$ws.Range("W25").Value = 1.49
$ws.Range("G26").Value = 3.6
$ws.Range("Z26").Value = 19
$ws.Range("AA26").Value = 40
$ws.Range("AE26").Value = 38
$ws.Range("AF26").Value = 27
$ws.Range("AH26").Value = 24
$ws.Range("AN26").Value = 55
$ws.Range("M27").Value = 1.05
$ws.Range("N27").Value = 2.7
$ws.Range("O27").Value = 1.06
$ws.Range("T27").Value = 1.04
$ws.Range("U27").Value = 1.04
$ws.Range("V27").Value = 1.14
$ws.Range("K29").Value = 3.8
$ws.Range("M29").Value = 1.05
$ws.Range("N29").Value = 1.1
$ws.Range("U29").Value = 2.02
$ws.Range("G30").Value = 2.08
$ws.Range("J30").Value = 3.55
$ws.Range("P30").Value = 1.78
$ws.Range("R30").Value = 1.3
$ws.Range("U30").Value = 1.94
$ws.Range("V30").Value = 1.29
$ws.Range("W30").Value = 1.93
$ws.Range("Y30").Value = 16.5
$ws.Range("AI30").Value = 80
$ws.Range("T31").Value = 1.04
$ws.Range("U31").Value = 1.04
$ws.Range("G32").Value = 2.16
$ws.Range("J32").Value = 2.76
$ws.Range("N32").Value = 2.42
$ws.Range("P32").Value = 2.14
$ws.Range("Q32").Value = 1.41
$ws.Range("T32").Value = 1.5
$ws.Range("U32").Value = 1.04
$ws.Range("W32").Value = 1.86
$ws.Range("N34").Value = 2.44
$ws.Range("T34").Value = 1.38
$ws.Range("U34").Value = 2.12
$ws.Range("F35").Value = 2.14
$ws.Range("L35").Value = 1.26
$ws.Range("K36").Value = 500
$ws.Range("N36").Value = 1.1
$ws.Range("P36").Value = 2.6
$ws.Range("S36").Value = 1.92
$ws.Range("N37").Value = 1.1
$ws.Range("P37").Value = 2.08
$ws.Range("R37").Value = 1.45
$ws.Range("S37").Value = 2.02
$ws.Range("T37").Value = 1.04
$ws.Range("U37").Value = 1.04
$ws.Range("N38").Value = 1.1
$ws.Range("T38").Value = 1.34
$ws.Range("U38").Value = 1.04
$ws.Range("T39").Value = 1.33
$ws.Range("U39").Value = 2.16
$ws.Range("V39").Value = 1.39
$ws.Range("Q40").Value = 1.29
$ws.Range("T40").Value = 1.04
$ws.Range("U40").Value = 1.04
$ws.Range("K42").Value = 3.25
$ws.Range("U42").Value = 1.62
$ws.Range("G43").Value = 2.68
$ws.Range("M43").Value = 1.06
$ws.Range("N43").Value = 2.24
$ws.Range("S43").Value = 3.05
$ws.Range("T43").Value = 1.04
$ws.Range("U43").Value = 1.04
$ws.Range("W43").Value = 1.6
$ws.Range("N44").Value = 2.6
$ws.Range("T44").Value = 1.04
$ws.Range("U44").Value = 1.04
$ws.Range("P45").Value = 2.14
$ws.Range("U45").Value = 2.28
$ws.Range("W45").Value = 1.33
$ws.Range("AM45").Value = 80
$ws.Range("F46").Value = 1.75
$ws.Range("U46").Value = 1.89
$ws.Range("K47").Value = 3.4
$ws.Range("P47").Value = 1.71
$ws.Range("N48").Value = 2.58
$ws.Range("O48").Value = 1.61
$ws.Range("Q48").Value = 2.8
$ws.Range("R48").Value = 1.2
$ws.Range("Y48").Value = 9.6
$ws.Range("F49").Value = 3.2
$ws.Range("N49").Value = 2.8
$ws.Range("T49").Value = 1.04
$ws.Range("U49").Value = 1.04
$ws.Range("N51").Value = 2.58
$ws.Range("T51").Value = 1.04
$ws.Range("U51").Value = 1.04
